$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.38

$ws.Range("B3").Value = 0.17
$ws.Range("C3").Value = 0.1

$ws.Range("B4").Value = 0.16
$ws.Range("C4").Value = 0.31

$ws.Range("B5").Value = 0.27
$ws.Range("C5").Value = 0.15

$ws.Range("A6").Value = "Surprise"
$ws.Range("B6").Value = 0.09
$ws.Range("C6").Value = 0.05

$ws.Rows(7).Delete()
